$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Shared-string text correction: "E7760" -> "E7420" (column K, rows 2-32
#    all shared the same text via the shared-string table).
$ws.Range("K2:K32").Value = "E7420"

# 2. Column K adopts the boolean-style number format that column L was
#    already using (style id 6 in the original workbook: numFmt 165 -
#    "TRUE";"TRUE";"FALSE").
$ws.Range("K2:K32").NumberFormat = $ws.Range("L2").NumberFormat

# 3. Column L gets a new font (size 11, family 0 instead of the inherited
#    size 10 / family 2) - apply once to L2, then propagate the resulting
#    cell format (not the named "Style") to the rest of the column so a
#    single new font/xf pair is reused by every row.
$l2 = $ws.Range("L2")
$l2.Font.Size = 11
$l2.Font.Family = 0
$l2.Copy()
$ws.Range("L2:L32").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 4. Column L switches from a cached literal boolean to a live formula
#    (=FALSE()) - write cell-by-cell so each row keeps its own <f> entry
#    instead of being collapsed into a shared-formula group.
for ($row = 2; $row -le 32; $row++) {
    $ws.Cells.Item($row, 12).Formula = "=FALSE()"
}

# 5. The active selection moves from L2:L32 to K2:K32.
$ws.Range("K2:K32").Select()
